$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.85%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "18"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.99%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "18"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.717"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.45%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "18"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08082"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.25%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "18"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.35%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "18"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.751"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.29%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "18"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.543"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.37%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "18"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.04%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "18"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9242"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.10%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "18"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1268"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.98%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "18"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1960"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.24%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "18"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.669"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12.39%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "18"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09194"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.99%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "18"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03706"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "7.76%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "18"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.12%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "18"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001298"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.90%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "18"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006222"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.46%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "18"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.376"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.10%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "18"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3485"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.31%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "18"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1421"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.89%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "18"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2659"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.63%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "18"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04437"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.44%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "18"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.13%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "18"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004325"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.20%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "18"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001242"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.40%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "18"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "18"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "18"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "18"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "18"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "18"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "18"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "18"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "18"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "18"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "18"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "18"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "18"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02895"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "16.46%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "18"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05484"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.72%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "18"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007827"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.85%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "18"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009936"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "11.09%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "18"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1420"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.20%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "18"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002216"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.31%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "18"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01105"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.10%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "18"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006801"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.34%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "18"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.38%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "18"

$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.90%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "18"

$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002287"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "27.10%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "18"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.38%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "18"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.38%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "18"

Write-Output "Applied updates to cryptos sheet"